$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tickers = @(
    "AAF",
    "ABDN",
    "ABF",
    "ANTO",
    "AUTO",
    "AV",
    "BARC",
    "BATS",
    "BDEV",
    "BEZ",
    "BF.B",
    "BKG",
    "BNZL",
    "BRBY",
    "BRK.B",
    "BT-A",
    "CCH",
    "CRDA",
    "DCC",
    "DGE",
    "ENT",
    "EXPN",
    "FCIT",
    "FRAS",
    "GLEN",
    "HLMA",
    "HSBA",
    "HSX",
    "IMB",
    "IMI",
    "INF",
    "ITRK",
    "JMAT",
    "KGF",
    "LGEN",
    "LLOY",
    "LSEG",
    "MNDI",
    "MNG",
    "OCDO",
    "PHNX",
    "PSON",
    "REL",
    "RMV",
    "RR",
    "RS1",
    "SBRY",
    "SDR",
    "SGRO",
    "SKG",
    "SMDS",
    "SMT",
    "SN",
    "SPX",
    "SSE",
    "STAN",
    "STJ",
    "ULVR",
    "UU",
    "WEIR",
    "WTB"
)

$startRow = 4979
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}
